$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = "O=[N+]([O-])c1cc(P(=O)(O)O)ccc1"
$ws.Range("E10").Value = "no_sanitizable"
$ws.Range("D11").Value = "Clc1c(Cl)c(Cl)c2oc3c(c(Cl)c(Cl)c(Cl)c3)c2c1"
$ws.Range("D13").Value = "Nc1cc2cc(S(=O)(=O)O)c(N=Nc3ccc(-c4ccc(N=Nc5c(S(=O)(=O)O)cc6cc(N)ccc6c5[O-])cc4)cc3)c([O-])c2cc1.[Na+].[Na+]"
